# define.xlsx: rework the "define-funcexec" demo sheet into a smaller
# "bak" demo — new header labels, a simple A2+B2 -> C2 -> D2 formula
# chain, and a single DEFINE() example (renamed from DEVDEF_1/DEVDEF_2
# to "c_from_a"); also drops the old Output1/Output2 + DEVDEF_2 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old demo layout (A1:G7: Input1/Input2/Constant/Multiplier/
# Formula/Output1/Output2 headers, row2 calc, row6/row7 DEFINE calls).
$ws.Range("A1:G7").Clear()

# New header row.
$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"
$ws.Range("C1").Value = "formula"
$ws.Range("D1").Value = "c(a)"

# New calc row: B2 is a literal input, C2/D2 chain off it.
$ws.Range("B2").Value = 5
$ws.Range("C2").Formula = "=A2+B2"
$ws.Range("D2").Formula = "=C2"

# Single DEFINE() demo row (replaces the old DEVDEF_1 / DEVDEF_2 rows).
$ws.Range("A6").Value = "Define 1"
$ws.Range("B6").Formula = "=DEFINE(""c_from_a"",A2,""#"",D2)"

# Selection moves off the old H7 anchor.
$null = $ws.Range("I11").Select()

# Window is minimized and nudged down; recalculation goes single-threaded.
$win = $excel.ActiveWindow
$win.WindowState = -4140
$win.Top = 900
$excel.MultiThreadedCalculation.Enabled = $false
